$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$xlPasteFormats = -4122

# ---- Phase 1: copy formats from existing styled cells to their new positions ----
# (ordered so every source cell is read before it is overwritten)
$ws.Range("A1").Copy() | Out-Null
$ws.Range("B1").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("A1").Copy() | Out-Null
$ws.Range("C1").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("A2").Copy() | Out-Null
$ws.Range("B2").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("B5").Copy() | Out-Null
$ws.Range("A3").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("C2").Copy() | Out-Null
$ws.Range("B3").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("A4").Copy() | Out-Null
$ws.Range("C3").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("B5").Copy() | Out-Null
$ws.Range("B4").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("A4").Copy() | Out-Null
$ws.Range("C4").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("A4").Copy() | Out-Null
$ws.Range("A5").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("A4").Copy() | Out-Null
$ws.Range("C5").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("A7").Copy() | Out-Null
$ws.Range("A6").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("A4").Copy() | Out-Null
$ws.Range("B6").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("A4").Copy() | Out-Null
$ws.Range("C6").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("B5").Copy() | Out-Null
$ws.Range("B7").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("A4").Copy() | Out-Null
$ws.Range("C7").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("A9").Copy() | Out-Null
$ws.Range("A8").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("A9").Copy() | Out-Null
$ws.Range("A10").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("A9").Copy() | Out-Null
$ws.Range("A11").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("A9").Copy() | Out-Null
$ws.Range("A12").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("A9").Copy() | Out-Null
$ws.Range("A13").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("A9").Copy() | Out-Null
$ws.Range("A18").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("B9").Copy() | Out-Null
$ws.Range("B18").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A19").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("C2").Copy() | Out-Null
$ws.Range("B19").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("C2").Copy() | Out-Null
$ws.Range("C19").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("A4").Copy() | Out-Null
$ws.Range("A7").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("A4").Copy() | Out-Null
$ws.Range("B5").PasteSpecial($xlPasteFormats) | Out-Null

$excel.CutCopyMode = $false

# ---- Phase 2: clear cells that no longer hold content ----
$ws.Range("B8").Clear() | Out-Null
$ws.Range("C8").Clear() | Out-Null
$ws.Range("B9").Clear() | Out-Null
$ws.Range("C9").Clear() | Out-Null

# ---- Phase 3: set cell values (and bold-run formatting for rich text) ----
$ws.Range("A1").Value = "Scenario"
$ws.Range("B1").Value = "Commands"
$ws.Range("C1").Value = "Description"
$ws.Range("A2").Value = "Case 1: Clone complete github repository to your system"
$ws.Range("B2").Value = "git clone https://github.com/drsladit/Github-LearningRepository.git"
$ws.Range("C2").Value = "1) If you want a copy of existing repository from github use git clone. git clone is how you get a local copy of an existing repository to work on.`n2) You no need to initialise the repository when you are using git clone i.e. No need to use command git init before typing git clone URL command`n3) It's usually only used once for a given repository, unless you want to have multiple working copies of it around. (Or want to get a clean copy after messing up your local one...)`n"
$ws.Range("A3").Value = "Case 2: Create new local repository and push changes to github repository`nInitialising: New local repository"
$ws.Range("A3").Characters(1, 88).Font.Bold = $true
$ws.Range("B3").Value = "First create a respository in Github account. Then clone the repository from github to your local system.`ngit init - if you do not have .git folder"
$ws.Range("C3").Value = "1) To initialize existing folder"
$ws.Range("A4").Value = "Staging: Adding files to staging "
$ws.Range("A4").Characters(1, 7).Font.Bold = $true
$ws.Range("B4").Value = "git add <Filename>      -- To add file into staging`ngit add *.py   -- Will add all .py files`ngit add .   -- Will add all files in the folder"
$ws.Range("C4").Value = "2) To add file into staging"
$ws.Range("A5").Value = "Staging status: To check files that are in staging"
$ws.Range("A5").Characters(1, 14).Font.Bold = $true
$ws.Range("B5").Value = "git status "
$ws.Range("A6").Value = "Staging removing files"
$ws.Range("B6").Value = "git rm --cached <FileName> -- To remove file from staging to untrack file"
$ws.Range("A7").Value = "Commit changes: From staging to local repository"
$ws.Range("A7").Characters(1, 14).Font.Bold = $true
$ws.Range("B7").Value = "git commit   --To commit changes to global repository with edit mode open`ngit commit -m 'changed GitCommands.txt' - To commit a file with out edit mode"
$ws.Range("A8").Value = "Stash"
$ws.Range("A9").Value = "merge"
$ws.Range("A10").Value = "commit - soft and hard reset"
$ws.Range("A11").Value = "checkout"
$ws.Range("A12").Value = "creating new branch"
$ws.Range("A13").Value = "How to navigate to existing branch"
$ws.Range("A18").Value = "Push changes: Push changes from local to Central/origin/github"
$ws.Range("A18").Characters(1, 13).Font.Bold = $true
$ws.Range("B18").Value = "git push - will be sufficient.`n"
$ws.Range("C18").Value = "1) After executing 2 commands you can see the changes commited to github repository"
$ws.Range("A19").Value = "Pull: To update local repository by pulling from central/remote"
$ws.Range("A19").Characters(1, 5).Font.Bold = $true
$ws.Range("B19").Value = "git init`ngit remote add origin git://github.com/cmcculloh/repo.git`ngit fetch --all`ngit pull origin master"
$ws.Range("C19").Value = "1) git pull (or git fetch + git merge) is how you update that local copy with new commits from the remote repository. If you are collaborating with others, it is a command that you will run frequently.`n`nAs your first example shows, it is possible to emulate git clone with an assortment of other git commands, but it's not really the case that git pull is doing `"basically the same thing`" as git clone (or vice-versa).`n"

# ---- Phase 4: new style for A2 (Case 1 heading) built on the bold+border style ----
$c = $ws.Range("A2")
$c.HorizontalAlignment = -4131
$c.VerticalAlignment = -4108
$c.WrapText = $true

# ---- Phase 5: row heights ----
$ws.Rows.Item(2).RowHeight = 120
$ws.Rows.Item(3).RowHeight = 45
$ws.Rows.Item(4).RowHeight = 45
$ws.Rows.Item(7).RowHeight = 60
$ws.Rows.Item(18).RowHeight = 30
$ws.Rows.Item(19).RowHeight = 120

# ---- Phase 6: selection ----
$ws.Range("B7").Select()
